$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Version: 1.0 -> 1.2.5
$ws.Range("D2").Value = "1.2.5"

# Precondition text correction (shared across all test cases in the sheet)
$ws.Range("B8").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B15").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B23").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B31").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B40").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B49").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B59").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B68").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B77").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B86").Value = "O usuário devidamente autenticado e na tela inicial do sistema."

# TC2 expected result wording fix
$ws.Range("D18").Value = "SYSTEM Identifica que a prestação de contas indicada pelo usuário não está em nenhum desses dois estados: a) NÃO REALIZADA e b) DEVOLVIDA; Não permite um novo envio ou alterações na prestação (exclusão de documentos)."

# TC4 expected result: add trailing period
$ws.Range("D35").Value = "SYSTEM Apresenta a tela de Detalhar Diárias."

# TC10 step: fix typo "histório" -> "histórico"
$ws.Range("B90").Value = "Chefe Verifica o histórico da tramitação da prestação de contas."
